$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.0101971532915357
$ws.Range("D2").Value = 0.003626923131477611
$ws.Range("E2").Value = 0.05289587296702258
$ws.Range("F2").Value = 0.4411102620360623
$ws.Range("G2").Value = 0.3200767524363641
$ws.Range("H2").Value = 0.3880983609024327
$ws.Range("I2").Value = 0.321948227514099
$ws.Range("M2").Value = 0.8391109518968847
$ws.Range("N2").Value = 1.806390725800441
$ws.Range("O2").Value = 1.343800827141138
$ws.Range("C3").Value = 0.009034222291994354
$ws.Range("D3").Value = 0.003166069604656485
$ws.Range("E3").Value = 0.05639259184957734
$ws.Range("F3").Value = 0.4172000692780244
$ws.Range("G3").Value = 0.2955380779540917
$ws.Range("H3").Value = 0.3805128844400656
$ws.Range("I3").Value = 0.3045092174425434
$ws.Range("M3").Value = 0.7373166087641039
$ws.Range("N3").Value = 1.651773543296031
$ws.Range("O3").Value = 1.275564322475958
$ws.Range("C4").Value = 0.008321623378478193
$ws.Range("D4").Value = 0.002882488985356702
$ws.Range("E4").Value = 0.05870294115950392
$ws.Range("F4").Value = 0.4028256908167691
$ws.Range("G4").Value = 0.2806897144269414
$ws.Range("H4").Value = 0.3761042544726649
$ws.Range("I4").Value = 0.2940277215863034
$ws.Range("M4").Value = 0.6746455644434803
$ws.Range("N4").Value = 1.556971036514909
$ws.Range("O4").Value = 1.23463823052407
$ws.Range("C5").Value = 0.008031608154425385
$ws.Range("D5").Value = 0.002766783468775458
$ws.Range("E5").Value = 0.05968510855047571
$ws.Range("F5").Value = 0.3970448364041914
$ws.Range("G5").Value = 0.2746934401400267
$ws.Range("H5").Value = 0.3743702073086581
$ws.Range("I5").Value = 0.2898130748586638
$ws.Range("M5").Value = 0.6490658801042599
$ws.Range("N5").Value = 1.518376361941876
$ws.Range("O5").Value = 1.218203792951044
$ws.Range("C6").Value = 0.007983474301582305
$ws.Range("D6").Value = 0.002747562260218217
$ws.Range("E6").Value = 0.05985064221806446
$ws.Range("F6").Value = 0.3960895619008653
$ws.Range("G6").Value = 0.27370104707785
$ws.Range("H6").Value = 0.3740860429972912
$ws.Range("I6").Value = 0.2891166511493779
$ws.Range("M6").Value = 0.6448159850434223
$ws.Range("N6").Value = 1.511970185676631
$ws.Range("O6").Value = 1.215489531610558
$ws.Range("C7").Value = 0.008317710599740735
$ws.Range("D7").Value = 0.002880929109828401
$ws.Range("E7").Value = 0.05871602281820243
$ws.Range("F7").Value = 0.4027474175036474
$ws.Range("G7").Value = 0.2806086261962832
$ws.Range("H7").Value = 0.3760806155301566
$ws.Range("I7").Value = 0.2939706523049423
$ws.Range("M7").Value = 0.6743007504402954
$ws.Range("N7").Value = 1.556450374696055
$ws.Range("O7").Value = 1.23441560624596
$ws.Range("C8").Value = 0.009795881040652432
$ws.Range("D8").Value = 0.003468154548283309
$ws.Range("E8").Value = 0.05406740150849587
$ws.Range("F8").Value = 0.43280222718505
$ws.Range("G8").Value = 0.3115702324801504
$ws.Range("H8").Value = 0.3854311821854992
$ws.Range("I8").Value = 0.315888189793526
$ws.Range("M8").Value = 0.8040483323266301
$ws.Range("N8").Value = 1.753053829675565
$ws.Range("O8").Value = 1.320070648128478
$ws.Range("C9").Value = 0.0127056890284436
$ws.Range("D9").Value = 0.004614418142605103
$ws.Range("E9").Value = 0.04626510282778673
$ws.Range("F9").Value = 0.4941878322768503
$ws.Range("G9").Value = 0.3740400185152168
$ws.Range("H9").Value = 0.4057474478859717
$ws.Range("I9").Value = 0.3606748108487778
$ws.Range("M9").Value = 1.05707856001996
$ws.Range("N9").Value = 2.139465154819618
$ws.Range("O9").Value = 1.495802234160237
$ws.Range("C10").Value = 0.01485004507391352
$ws.Range("D10").Value = 0.005452888770324904
$ws.Range("E10").Value = 0.04135746289745112
$ws.Range("F10").Value = 0.540807067454125
$ws.Range("G10").Value = 0.4210391953494366
$ws.Range("H10").Value = 0.4218899593114145
$ws.Range("I10").Value = 0.3947013060638369
$ws.Range("M10").Value = 1.242056010451932
$ws.Range("N10").Value = 2.423678289929512
$ws.Range("O10").Value = 1.629734618383566
$ws.Range("C11").Value = 0.01582694034210164
$ws.Range("D11").Value = 0.005833434950787364
$ws.Range("E11").Value = 0.03930907181283283
$ws.Range("F11").Value = 0.5623512875463348
$ws.Range("G11").Value = 0.4426671690091553
$ws.Range("H11").Value = 0.4294997986830253
$ws.Range("I11").Value = 0.4104290466843992
$ws.Range("M11").Value = 1.325993434783641
$ws.Range("N11").Value = 2.553000083283166
$ws.Range("O11").Value = 1.691731319612359
$ws.Range("C12").Value = 0.01619706012970568
$ws.Range("D12").Value = 0.005977401857197862
$ws.Range("E12").Value = 0.03856030642766139
$ws.Range("F12").Value = 0.570558341765377
$ws.Range("G12").Value = 0.4508932444487925
$ws.Range("H12").Value = 0.4324198992350716
$ws.Range("I12").Value = 0.4164208175271114
$ws.Range("M12").Value = 1.357746718811839
$ws.Range("N12").Value = 2.601971329837852
$ws.Range("O12").Value = 1.715363059083643
$ws.Range("C13").Value = 0.01611733993573949
$ws.Range("D13").Value = 0.005946402311831633
$ws.Range("E13").Value = 0.0387203623341601
$ws.Range("F13").Value = 0.568788634254588
$ws.Range("G13").Value = 0.4491200031250457
$ws.Range("H13").Value = 0.4317892931329652
$ws.Range("I13").Value = 0.4151287770704499
$ws.Range("M13").Value = 1.350909533074955
$ws.Range("N13").Value = 2.591424574538905
$ws.Range("O13").Value = 1.710266633836
$ws.Range("C14").Value = 0.01585738657174574
$ws.Range("D14").Value = 0.005845282019482312
$ws.Range("E14").Value = 0.03924692847323863
$ws.Range("F14").Value = 0.5630255086991838
$ws.Range("G14").Value = 0.4433432088022471
$ws.Range("H14").Value = 0.4297392665960444
$ws.Range("I14").Value = 0.4109212701731906
$ws.Range("M14").Value = 1.3286064482905
$ws.Range("N14").Value = 2.557029006038022
$ws.Range("O14").Value = 1.693672405963525
$ws.Range("C15").Value = 0.01569818208430718
$ws.Range("D15").Value = 0.005783324653233279
$ws.Range("E15").Value = 0.03957298376665097
$ws.Range("F15").Value = 0.5595017817120151
$ws.Range("G15").Value = 0.4398094595427438
$ws.Range("H15").Value = 0.428488572265735
$ws.Range("I15").Value = 0.4083487477630712
$ws.Range("M15").Value = 1.314940943742471
$ws.Range("N15").Value = 2.535960574832643
$ws.Range("O15").Value = 1.683528170163413
$ws.Range("C16").Value = 0.01478623007326263
$ws.Range("D16").Value = 0.00542800053651149
$ws.Range("E16").Value = 0.04149507256093976
$ws.Range("F16").Value = 0.5394059052049442
$ws.Range("G16").Value = 0.4196307769062457
$ws.Range("H16").Value = 0.4213980096074863
$ws.Range("I16").Value = 0.393678489137784
$ws.Range("M16").Value = 1.236566124437871
$ws.Range("N16").Value = 2.415227044953269
$ws.Range("O16").Value = 1.625704612226286
$ws.Range("C17").Value = 0.01422713074018844
$ws.Range("D17").Value = 0.005209787647054043
$ws.Range("E17").Value = 0.04272168129706655
$ws.Range("F17").Value = 0.5271642228601365
$ws.Range("G17").Value = 0.4073155615326414
$ws.Range("H17").Value = 0.4171165097942975
$ws.Range("I17").Value = 0.3847426730520738
$ws.Range("M17").Value = 1.188430678893397
$ws.Range("N17").Value = 2.34116573104825
$ws.Range("O17").Value = 1.590506575440088
$ws.Range("C18").Value = 0.01390568574429096
$ws.Range("D18").Value = 0.005084195372841549
$ws.Range("E18").Value = 0.04344449791859528
$ws.Range("F18").Value = 0.5201548245143215
$ws.Range("G18").Value = 0.4002554770368079
$ws.Range("H18").Value = 0.4146789858317277
$ws.Range("I18").Value = 0.3796264439733363
$ws.Range("M18").Value = 1.160724810156836
$ws.Range("N18").Value = 2.298570892747136
$ws.Range("O18").Value = 1.57036226455898
$ws.Range("C19").Value = 0.01379687342704727
$ws.Range("D19").Value = 0.005041658250352299
$ws.Range("E19").Value = 0.04369218977357159
$ws.Range("F19").Value = 0.5177869975149321
$ws.Range("G19").Value = 0.3978690437335501
$ws.Range("H19").Value = 0.4138579870474928
$ws.Range("I19").Value = 0.3778981916374136
$ws.Range("M19").Value = 1.151340763179277
$ws.Range("N19").Value = 2.284149717821265
$ws.Range("O19").Value = 1.563559005219446
$ws.Range("C20").Value = 0.01428663402874264
$ws.Range("D20").Value = 0.005233025350371179
$ws.Range("E20").Value = 0.04258931265320731
$ws.Range("F20").Value = 0.5284640886451513
$ws.Range("G20").Value = 0.4086241227267777
$ws.Range("H20").Value = 0.4175696862448177
$ws.Range("I20").Value = 0.3856914808457432
$ws.Range("M20").Value = 1.193556824602581
$ws.Range("N20").Value = 2.349049379965891
$ws.Range("O20").Value = 1.594243040128219
$ws.Range("C21").Value = 0.01593373610271698
$ws.Range("D21").Value = 0.005874987327640468
$ws.Range("E21").Value = 0.03909152922949533
$ws.Range("F21").Value = 0.564716953745048
$ws.Range("G21").Value = 0.4450390122981958
$ws.Range("H21").Value = 0.4303403656464582
$ws.Range("I21").Value = 0.4121561376818477
$ws.Range("M21").Value = 1.335158287857965
$ws.Range("N21").Value = 2.567131858293749
$ws.Range("O21").Value = 1.698542319768933
$ws.Range("C22").Value = 0.01701131827241653
$ws.Range("D22").Value = 0.006293738854363795
$ws.Range("E22").Value = 0.03696261744795981
$ws.Range("F22").Value = 0.5886944366025659
$ws.Range("G22").Value = 0.4690484482617876
$ws.Range("H22").Value = 0.438910724612839
$ws.Range("I22").Value = 0.4296623572907095
$ws.Range("M22").Value = 1.427515687671374
$ws.Range("N22").Value = 2.709658591458378
$ws.Range("O22").Value = 1.767611486227111
$ws.Range("C23").Value = 0.01643609553017455
$ws.Range("D23").Value = 0.006070320946072627
$ws.Range("E23").Value = 0.03808433794246957
$ws.Range("F23").Value = 0.5758711115273059
$ws.Range("G23").Value = 0.4562147964454084
$ws.Range("H23").Value = 0.4343160376557194
$ws.Range("I23").Value = 0.4202996666003571
$ws.Range("M23").Value = 1.378240549150718
$ws.Range("N23").Value = 2.63359117803833
$ws.Range("O23").Value = 1.730664952567622
$ws.Range("C24").Value = 0.01425973262623614
$ws.Range("D24").Value = 0.005222520015312426
$ws.Range("E24").Value = 0.04264910165441194
$ws.Range("F24").Value = 0.5278763305445722
$ws.Range("G24").Value = 0.4080324596424418
$ws.Range("H24").Value = 0.4173647305186563
$ws.Range("I24").Value = 0.3852624591631866
$ws.Range("M24").Value = 1.191239394168647
$ws.Range("N24").Value = 2.345485231772273
$ws.Range("O24").Value = 1.592553499473468
$ws.Range("C25").Value = 0.0119173448579204
$ws.Range("D25").Value = 0.004304937552422672
$ws.Range("E25").Value = 0.04823273670708228
$ws.Range("F25").Value = 0.477316535216616
$ws.Range("G25").Value = 0.3569490515952509
$ws.Range("H25").Value = 0.4000384560281844
$ws.Range("I25").Value = 0.3483633313797014
$ws.Range("M25").Value = 0.9887845913642366
$ws.Range("N25").Value = 2.034859715736161
$ws.Range("O25").Value = 1.710266633836

Write-Output "applied 380 kV case updates"